# Update the student_import workbook with the new roster (Zack Milele /
# Erustus Baraza) in place of the old sample rows (Mapenzi Karani / James
# Baya), dropping the phone/email columns in favour of address + roll
# code, and re-pointing the selection / page orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- name ---
$ws.Range("A2").Value = "Zack Milele"
$ws.Range("A3").Value = "Erustus Baraza"

# --- birthday (store as a real date, custom yyyy-mm-dd display format) ---
$ws.Range("B2").Value = 43020
$ws.Range("B3").Value = 42056
$ws.Range("B2:B3").NumberFormat = "yyyy-mm-dd"

# --- sex ---
$ws.Range("C2").Value = "male"
$ws.Range("C3").Value = "male"

# --- address (replaces the old numeric value in this column) ---
$ws.Range("D2").Value = "Nairobi"
$ws.Range("D3").Value = "Nairobi"

# --- phone column is no longer populated ---
$ws.Range("E2:E3").ClearContents()

# --- email column: drop the mailto hyperlinks and leave the cells blank ---
$ws.Hyperlinks.Delete()
$ws.Range("F2:F3").ClearContents()

# --- roll (now an admission-style code instead of a raw number) ---
$ws.Range("G2").Value = "KPV/2019/769"
$ws.Range("G3").Value = "KPV/2019/775"

# --- selection + print orientation ---
$ws.Range("B3").Select() | Out-Null
$ws.PageSetup.Orientation = 1
